# Treasure Hunt participants.xlsx
# Insert a new "Phone Number" column between "College" and "Purchase_Type",
# shifting Purchase_Type/Payment_Mode one column to the right, and fill in
# the phone number for the existing participant row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift column D (and everything right of it) one column over so the new
# "Phone Number" column lands at D, pushing Purchase_Type -> E and
# Payment_Mode -> F.
$ws.Columns.Item(4).Insert()

$ws.Range("D1").Value = "Phone Number"

# The phone number is digits-only, so a plain .Value assignment would be
# stored as a number. Write it as a text formula, then paste-special just
# the value back over itself so it ends up a normal static text cell
# (matching how the phone number is authored as text, not a number).
$ws.Range("D2").Formula = '="9429510862"'
$ws.Range("D2").Copy()
$ws.Range("D2").PasteSpecial(-4163)
